# Reorders the comma-separated "Recorded By" values in column G so that
# the literal token "System" (case-sensitive) is moved to the front of the
# list, with the remaining tokens kept in their original relative order
# (i.e. the leading run up to and including the first "System" token is
# rotated to the end). Rows whose list already starts with "System", rows
# without a "System" token at all, and single-value rows are left as-is.

function CaseSensitiveEquals($s1, $s2) {
    if ($s1.Length -ne $s2.Length) { return $false }
    for ($i = 0; $i -lt $s1.Length; $i++) {
        if ([int][char]$s1[$i] -ne [int][char]$s2[$i]) { return $false }
    }
    return $true
}

function Transform-RecordedBy($s) {
    $parts = $s.Split(",")
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    if ($trimmed.Count -le 1) { return $s }

    $hasSystem = $false
    foreach ($p in $trimmed) {
        if (CaseSensitiveEquals $p "System") { $hasSystem = $true }
    }
    if (-not $hasSystem) { return $s }

    $firstIsSystem = CaseSensitiveEquals $trimmed[0] "System"
    if ($firstIsSystem) { return $s }

    # Rotate left by one: move the first element to the end, keep the rest
    # in their original relative order.
    $rotated = $trimmed[1..($trimmed.Count - 1)] + $trimmed[0]
    return ($rotated -join ", ")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Column G = "Recorded By" (7th column); skip the header row.
$col = 7

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2
    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val -eq "Recorded By") { continue }

    $newVal = Transform-RecordedBy $val
    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
